$p = $ppt.ActivePresentation

# --- 1. Delete the old "Vi du:" slide (slide 2). All following slides
#        shift up by one position; their content (incl. a16/p14 creationIds)
#        stays attached to their own XML parts, which is exactly what the
#        target revision shows (each slide's content appears to "shift"
#        because the in-between slide was removed). ---------------------
$p.Slides.Item(2).Delete()

# --- 2. On the (now last, position 8) "Bai tap tu lam" slide, fix the
#        wording of Cau 2a: the pen colour changes from "mau xanh" (blue)
#        to "mau do" (red). The target text lives in a single run, so we
#        replace that whole run's text in one shot to avoid splitting it
#        into multiple runs. -------------------------------------------
$lastSlide = $p.Slides.Item($p.Slides.Count)
$needle = ", nét vẽ 2, bút vẽ màu xanh như Hình số 2."
$replacement = ", nét vẽ 2, bút vẽ màu đỏ như Hình số 2."
foreach ($shp in $lastSlide.Shapes) {
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $idx = $tr.Text.IndexOf($needle)
        if ($idx -ge 0) {
            $run = $tr.Characters($idx + 1, $needle.Length)
            $run.Text = $replacement
        }
    }
}
